# Collection_QRS_EQ5D-5L: add two blank "derivation" columns
# (derived_variable / derivation_description) right before the existing
# "codelist" column, i.e. insert two new columns at X:Y.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Collection_QRS_EQ5D-5L")

# Insert two blank columns at X:Y - this shifts the old X:AH block
# (codelist ... change_history) two columns to the right, to Z:AJ,
# carrying cell values/styles and column widths along with it.
$ws.Columns("X:Y").Insert()

# New columns inherit default width; make them match the width of the
# neighboring "display_hidden" column (W) like the rest of the sheet.
$ws.Columns("X:Y").ColumnWidth = $ws.Columns("W").ColumnWidth

# Label the two new (still blank) header cells.
$ws.Range("X1").Value = "derived_variable"
$ws.Range("Y1").Value = "derivation_description"

# Re-establish the autofilter / filter-database range over the now-wider
# used range (A1:AJ16 instead of A1:AH16).
$ws.AutoFilterMode = $false
$ws.Range("A1:AJ16").AutoFilter()

$filterDb = $wb.Names.Item("Collection_QRS_EQ5D-5L!_FilterDatabase")
$filterDb.RefersTo = "='Collection_QRS_EQ5D-5L'!`$A`$1:`$AJ`$16"

# Restore/update the view: scroll the frozen header row over a bit and
# move the active cell in the lower pane onto the new columns.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 14
$ws.Range("X4").Select()
